$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "MuSCs" data row (row 3) entirely - this also shifts the
# shared string table so "Avp"/"Oxtr" indices shift down and "MuSCs" is dropped.
$ws.Rows.Item(3).Delete()

# Refresh a handful of numeric values in row 2 to reflect the updated TPM-based
# calculations from the commit.
$ws.Range("G2").Value = 0.3143816666666667
$ws.Range("H2").Value = 0.943145
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.05993299999999999
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
